$wb = $excel.ActiveWorkbook

# --- Fix survey_1 (sheet "survey_1") rows 15-17: restore names to match column A ---
$survey1 = $wb.Worksheets.Item("survey_1")
$survey1.Range("B15").Value = $survey1.Range("A15").Value2
$survey1.Range("B16").Value = $survey1.Range("A16").Value2
$survey1.Range("B17").Value = $survey1.Range("A17").Value2

# --- Fix survey_2 (sheet "survey_2") rows 5-7: restore names to match column A ---
$survey2 = $wb.Worksheets.Item("survey_2")
$survey2.Range("B5").Value = $survey2.Range("A5").Value2
$survey2.Range("B6").Value = $survey2.Range("A6").Value2
$survey2.Range("B7").Value = $survey2.Range("A7").Value2

# --- Fix survey_2 rows 188, 192, 196: older_adult contact variable names ---
$survey2.Range("B188").Value = "multiple_contacts_older_adult_work"
$survey2.Range("B192").Value = "multiple_contacts_older_adult_school"
$survey2.Range("B196").Value = "multiple_contacts_older_adult_other"

# --- Insert a new "locations" worksheet right before "oldsheet_v2" ---
$oldsheet = $wb.Worksheets.Item("oldsheet_v2")
$locations = $wb.Worksheets.Add($oldsheet)
$locations.Name = "locations"

$locData = @(
    ,@('country','variable','oldname','newname')
    ,@('be','area_3_name','Région de Bruxelles-Capitale / Brussels Hoofdsted','Région de Bruxelles-Capitale / Brussels Hoofdstede')
    ,@('be','area_3_name','Région de Bruxelles-Capitale / Brussels Hoofdstede','Région de Bruxelles-Capitale / Brussels Hoofdstede')
    ,@('be','area_3_name','Région Flamande / Vlaams Gewest','Région Flamande / Vlaams Gewest')
    ,@('be','area_3_name','Région Wallonne / Waals Gewest','Région Wallonne / Waals Gewest')
    ,@('nl','area_3_name','Noord-Nederland','Noord-Nederland')
    ,@('nl','area_3_name','Oost-Nederland','Oost-Nederland')
    ,@('nl','area_3_name','West-Nederland','West-Nederland')
    ,@('nl','area_3_name','Zuid-Nederland','Zuid-Nederland')
    ,@('no','area_3_name','Akershus','Akershus')
    ,@('no','area_3_name','Aust-Agder','Aust-Agder')
    ,@('no','area_3_name','Buskerud','Buskerud')
    ,@('no','area_3_name','Finnmark','Finnmark')
    ,@('no','area_3_name','Hedmark','Hedmark')
    ,@('no','area_3_name','Hordaland','Hordaland')
    ,@('no','area_3_name','Møre og Romsdal','Møre og Romsdal')
    ,@('no','area_3_name','Nordland','Nordland')
    ,@('no','area_3_name','Oppland','Oppland')
    ,@('no','area_3_name','Oslo','Oslo')
    ,@('no','area_3_name','Østfold','Østfold')
    ,@('no','area_3_name','Rogaland','Rogaland')
    ,@('no','area_3_name','Sogn og Fjordane','Sogn og Fjordane')
    ,@('no','area_3_name','Telemark','Telemark')
    ,@('no','area_3_name','Troms','Troms')
    ,@('no','area_3_name','Trøndelag','Trøndelag')
    ,@('no','area_3_name','Vest-Agder','Vest-Agder')
    ,@('no','area_3_name','Vestfold','Vestfold')
    ,@('uk','area_3_name','East Anglia','East of England')
    ,@('uk','area_3_name','East Midlands','East Midlands')
    ,@('uk','area_3_name','East of Engla','East of England')
    ,@('uk','area_3_name','East of England','East of England')
    ,@('uk','area_3_name','Greater Londo','Greater London')
    ,@('uk','area_3_name','Greater London','Greater London')
    ,@('uk','area_3_name','North East','North East')
    ,@('uk','area_3_name','North West','North West')
    ,@('uk','area_3_name','Northern Irel','Northern Ireland')
    ,@('uk','area_3_name','Northern Ireland','Northern Ireland')
    ,@('uk','area_3_name','Scotland','Scotland')
    ,@('uk','area_3_name','South East','South East')
    ,@('uk','area_3_name','South West','South West')
    ,@('uk','area_3_name','Wales','Wales')
    ,@('uk','area_3_name','West Midlands','West Midlands')
    ,@('uk','area_3_name','Yorkshire and','Yorkshire and The Humber')
    ,@('uk','area_3_name','Yorkshire and Humberside','Yorkshire and The Humber')
    ,@('uk','area_3_name','Yorkshire and The Humber','Yorkshire and The Humber')
)

$rowCount = $locData.Count
$colCount = $locData[0].Count
$arr = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r,$c] = $locData[$r][$c]
    }
}

$targetRange = $locations.Range($locations.Cells.Item(1,1), $locations.Cells.Item($rowCount, $colCount))
$targetRange.Value = $arr

$locations.Columns.Item(1).ColumnWidth = 41.66
$locations.Columns.Item(2).ColumnWidth = 41.66

$locations.Range("D4:D5").Select()

# --- Restore view/selection state ---
$survey1.Activate()
$survey1.Range("B15:B17").Select()

$survey2.Activate()
$survey2.Range("B197").Select()
